# EI Variable Installments T2 scenarios
# Insert a new row (waittopageload1 / 2000) above the existing
# "clickonsubmit / Submit" row on the "Edit Repayment Schedule" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 6, shifting existing rows 6-12 down to 7-13.
$ws.Rows("6:6").Insert()

# Copy the formatting of row 3 (label/value pair styling) into the new row
# so the new cells pick up the same look (label style + numeric value style)
# as the other label/amount rows on this sheet.
$ws.Range("A3:B3").Copy($ws.Range("A6:B6"))

# Set the new row's content.
$ws.Range("A6").Value = "waittopageload1"
$ws.Range("B6").Value = 2000

# Match the workbook's recorded selection after the edit.
$ws.Range("A6:B6").Select() | Out-Null
